$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (style) from H1 (existing header) onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-57
$I_vals = @(7,5,8,8,8,5,4,7,6,8,8,8,7,7,6,6,6,8,9,8,6,4,8,9,6,7,8,9,9,8,4,6,8,7,6,5,10,5,6,8,3,6,6,7,7,5,7,5,9,7,8,7,6,6,6,3)
$J_vals = @(8,6,8,8,8,5,5,7,7,8,8,8,7,7,6,6,7,8,9,9,6,5,8,9,7,8,9,9,9,8,5,7,8,8,6,5,10,5,6,9,4,7,7,7,7,6,7,6,9,8,8,7,6,6,6,3)

for ($idx = 0; $idx -lt $I_vals.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $I_vals[$idx]
    $ws.Cells.Item($row, 10).Value = $J_vals[$idx]
}
